$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column ("11-dec") before the
# existing "01-oct." column (currently EN), shifting EN:FR right to
# EO:FS, and fill the new column's header + data cells. ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("EN1").EntireColumn.Insert()

$wsPrix.Range("EN1").Value = "11-dec"
$wsPrix.Range("EN2:EN25").Value = "-"

# --- Sheet "Gaz": append a new row with the latest data point. ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A174").Value = "'2025-12-09"
$wsGaz.Range("A174").Style = "Normal"
$wsGaz.Range("B174").Value = 26.1

# --- Sheet "CO2": append a new row with the latest data point. ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A174").Value = "'2025-12-09"
$wsCo2.Range("A174").Style = "Normal"
$wsCo2.Range("B174").Value = 82.67
